$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old units row (row 2); data rows 3-7 shift up to become rows 2-6.
$ws.Rows.Item(2).Delete()

# Build a temporary named style for the new header row: Arial 9, General
# format, applyFont only (matches the workbook's existing "data" font but
# without forcing a number format) - mirrors the xf that Excel produces when
# a style is assigned then detached, leaving only the cellXf behind.
$hdrStyle = $wb.Styles.Add("HeaderStyle")
$hdrStyle.Font.Name = "Arial"
$hdrStyle.Font.Size = 9

# Overwrite row 1 with the new header labels.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("A1:E1").ClearFormats()
$ws.Range("F1:K1").Style = "HeaderStyle"

# Detach the helper named style so the workbook's style tables stay exactly
# as they were (only the cellXf used by F1:K1 remains behind).
$wb.Styles.Item("HeaderStyle").Delete()

$ws.Range("A2:K2").Select()
